# JPS_1_AST/misc/Zapytania testowe JPS.xlsx
# "zaimplementowane asserty we wszystkich nowych testcase'ach"
#
# The author marked every B-column "expected result" cell belonging to a
# test case whose assert is now implemented with the sheet's existing
# green highlight style, un-hid the row that held the next case to look
# at, and left the window scrolled/zoomed/selected where they stopped
# working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Highlight the B-cells that were still plain (no fill) by copying
#    the format that is already used elsewhere in the column (green
#    fill + wrap, no special alignment). Using copy/paste-special keeps
#    reusing the workbook's existing fill definition instead of
#    creating a duplicate one. The target range is not contiguous, and
#    PasteSpecial only reliably applies to one contiguous block at a
#    time, so each block is pasted individually.
# ---------------------------------------------------------------------
$greenWrap = $ws.Range("B2")
$greenWrap.Copy()
$ws.Range("B56:B57").PasteSpecial(-4122)
$ws.Range("B60:B100").PasteSpecial(-4122)
$ws.Range("B102:B103").PasteSpecial(-4122)
$ws.Range("B105:B132").PasteSpecial(-4122)

# B101 gets the same green fill too, but it keeps its left alignment.
$ws.Range("B101").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B101").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B101").WrapText = $true

# ---------------------------------------------------------------------
# 2. The test row that used to be hidden (row 9) is shown again.
# ---------------------------------------------------------------------
$ws.Rows.Item(9).Hidden = $false

# ---------------------------------------------------------------------
# 3. Leave the window the way the author left it: zoomed to 115%,
#    scrolled further down, with A9 selected.
# ---------------------------------------------------------------------
$ws.Range("A9").Select()
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.Zoom = 115
